$d = $word.ActiveDocument

# Locate the "Git commit" paragraph (the one with the _GoBack bookmark).
$commitIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "Git commit*") {
        $commitIndex = $i
        break
    }
}

# Make room for a new paragraph right before it ("add  <nom de docs>").
$d.Paragraphs($commitIndex).Range.InsertParagraphBefore()

# After the insert, the blank new paragraph takes the original index and
# the "Git commit" paragraph shifts one slot later.
$addPara = $d.Paragraphs($commitIndex).Range
$commitRange = $d.Paragraphs($commitIndex + 1).Range

# Fill in the new "add  <nom de docs>" paragraph with the same run/proofErr
# layout the original "Git commit" paragraph used for its "commit" word.
$addXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:jc w:val="both"/>
            </w:pPr>
            <w:proofErr w:type="spellStart"/>
            <w:r>
              <w:t>add</w:t>
            </w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r>
              <w:t xml:space="preserve">  &lt;nom de docs&gt;</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@
$addPara.InsertXML($addXml)

# Rebuild the "Git commit" paragraph so "commit" gets its own spellEnd/run
# and is followed by a trailing space run before the existing bookmark.
# (Original paragraph-level rsid/paraId attributes are kept so this still
# reads as an in-place edit of the same paragraph.)
$commitXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
        <w:body>
          <w:p w14:paraId="54D2E278" w14:textId="6C3CE560" w:rsidR="005D6DD2" w:rsidRPr="005D6DD2" w:rsidRDefault="005D6DD2" w:rsidP="00DD5CB2">
            <w:pPr>
              <w:jc w:val="both"/>
            </w:pPr>
            <w:r w:rsidRPr="005D6DD2">
              <w:t xml:space="preserve">Git </w:t>
            </w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r w:rsidRPr="005D6DD2">
              <w:t>commit</w:t>
            </w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r>
              <w:t xml:space="preserve"> </w:t>
            </w:r>
            <w:bookmarkStart w:id="0" w:name="_GoBack"/>
            <w:bookmarkEnd w:id="0"/>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@
$commitRange.InsertXML($commitXml)
